$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 111
$ws.Range("H111").Value = 4
$ws.Range("I111").Value = 1
$ws.Range("J111").Value = "H"
$ws.Range("N111").Value = 1.45
$ws.Range("O111").Value = 4.2
$ws.Range("P111").Value = 5.5
$ws.Range("Q111").Value = -1
$ws.Range("S111").Value = 1.975
$ws.Range("U111").Value = 2
$ws.Range("V111").Value = 1.85
$ws.Range("W111").Value = 0.45
$ws.Range("X111").Value = -1
$ws.Range("Y111").Value = -1
$ws.Range("Z111").Value = 0.875
$ws.Range("AA111").Value = -1
$ws.Range("AB111").Value = 1
$ws.Range("AC111").Value = -1

# Row 112
$ws.Range("H112").Value = 3
$ws.Range("I112").Value = 4
$ws.Range("J112").Value = "A"
$ws.Range("N112").Value = 3.3
$ws.Range("O112").Value = 3.5
$ws.Range("P112").Value = 2
$ws.Range("R112").Value = 2.05
$ws.Range("S112").Value = 1.8
$ws.Range("T112").Value = 2.75
$ws.Range("U112").Value = 1.95
$ws.Range("V112").Value = 1.9
$ws.Range("W112").Value = -1
$ws.Range("X112").Value = -1
$ws.Range("Y112").Value = 1
$ws.Range("Z112").Value = -1
$ws.Range("AA112").Value = 0.8
$ws.Range("AB112").Value = 0.95
$ws.Range("AC112").Value = -1

# Row 113
$ws.Range("H113").Value = 2
$ws.Range("I113").Value = 1
$ws.Range("J113").Value = "H"
$ws.Range("N113").Value = 1.6
$ws.Range("O113").Value = 3.6
$ws.Range("P113").Value = 5.25
$ws.Range("R113").Value = 1.8
$ws.Range("S113").Value = 2.05
$ws.Range("U113").Value = 1.875
$ws.Range("W113").Value = 0.6000000000000001
$ws.Range("X113").Value = -1
$ws.Range("Y113").Value = -1
$ws.Range("Z113").Value = 0.4
$ws.Range("AA113").Value = -0.5
$ws.Range("AB113").Value = 0.875
$ws.Range("AC113").Value = -1
